$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "gray" sheet: append a new color-palette table (rows 17-27,
#    columns M:O) with a header row "dark / medium / light" followed
#    by ten rows of hex swatches (tab10 + tab10-medium + tab10-light).
# ------------------------------------------------------------------
$gray = $wb.Worksheets.Item("gray")
$gray.Activate() | Out-Null

$table = @(
    @("dark",    "medium",  "light"),
    @("#1F77B4", "#729ECE", "#AEC7E8"),
    @("#FF7F0E", "#FF9E4A", "#FFBB78"),
    @("#2CA02C", "#67BF5C", "#98DF8A"),
    @("#D62728", "#ED665D", "#FF9896"),
    @("#9467BD", "#AD8BC9", "#C5B0D5"),
    @("#8C564B", "#A8786E", "#C49C94"),
    @("#E377C2", "#ED97CA", "#F7B6D2"),
    @("#7F7F7F", "#A2A2A2", "#C7C7C7"),
    @("#BCBD22", "#CDCC5D", "#DBDB8D"),
    @("#17BECF", "#6DCCDA", "#9EDAE5")
)

$startRow = 17
for ($i = 0; $i -lt $table.Length; $i++) {
    $r = $startRow + $i
    $gray.Cells.Item($r, 13).Value = $table[$i][0]   # column M
    $gray.Cells.Item($r, 14).Value = $table[$i][1]   # column N
    $gray.Cells.Item($r, 15).Value = $table[$i][2]   # column O
}

# Leave the selection on the newly-added block, matching the edited
# workbook's saved cursor position.
$gray.Range("M18:O27").Select() | Out-Null

# ------------------------------------------------------------------
# 2) Move the active tab from "gray" to "light" (the last-active
#    sheet when the file was saved), updating tabSelected/activeTab.
# ------------------------------------------------------------------
$light = $wb.Worksheets.Item("light")
$light.Activate() | Out-Null
$light.Range("M13").Select() | Out-Null

Write-Output "edit complete"
